$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating the "2022-Q1" sheet
#    (same column layout / styles) and placing it right after "总计".
# ------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q1Template = $wb.Worksheets.Item("2022-Q1")
$q1Template.Copy($null, $zongji)

$q3 = $wb.Worksheets.Item("总计").Next
$q3.Name = "2022-Q3"

# The template only had 2 fund rows (rows 2-3); we need a 3rd (row 4).
# Clone row 3's formatting down to row 4 before filling in values.
$q3.Range("A3:H3").Copy()
$q3.Range("A4:H4").PasteSpecial(-4122)

# Row 2 -- new fund added this quarter.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'006282"
$q3.Range("C2").Value = "上投摩根欧洲动力策略股票（QDII）"
$q3.Range("D2").Value = "'0.41"
$q3.Range("E2").Value = "'91.47"
$q3.Range("F2").Value = "'1.97"
$q3.Range("G2").Value = "'0.0081"
$q3.Range("H2").Value = 8

# Row 3 -- 010343, updated figures.
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'010343"
$q3.Range("C3").Value = "华宝英国富时100指数（QDII）A"
$q3.Range("D3").Value = "'0.13"
$q3.Range("E3").Value = "'92.85"
$q3.Range("F3").Value = "'4.51"
$q3.Range("G3").Value = "'0.0059"
$q3.Range("H3").Value = 5

# Row 4 -- 010344, updated figures.
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'010344"
$q3.Range("C4").Value = "华宝英国富时100指数（QDII）C"
$q3.Range("D4").Value = "'0.08"
$q3.Range("E4").Value = "'92.85"
$q3.Range("F4").Value = "'4.51"
$q3.Range("G4").Value = "'0.0036"
$q3.Range("H4").Value = 5

# ------------------------------------------------------------------
# 2. Insert a new summary row for 2022-Q3 at the top of "总计"'s data
#    (row 2), pushing the existing quarters down by one row.
# ------------------------------------------------------------------
# Give the brand-new row 6 the same cell formatting as row 5 before
# shuffling values around.
$zongji.Range("A5").Copy()
$zongji.Range("A6").PasteSpecial(-4122)

$zongji.Range("A6").Value = 4
$zongji.Range("B6").Value = "2021-Q2"
$zongji.Range("C6").Value = 2
$zongji.Range("D6").Value = 0.01

$zongji.Range("A5").Value = 3
$zongji.Range("B5").Value = "2021-Q3"
$zongji.Range("C5").Value = 2
$zongji.Range("D5").Value = 0.01

$zongji.Range("A4").Value = 2
$zongji.Range("B4").Value = "2021-Q4"
$zongji.Range("C4").Value = 2
$zongji.Range("D4").Value = 0.01

$zongji.Range("A3").Value = 1
$zongji.Range("B3").Value = "2022-Q1"
$zongji.Range("C3").Value = 2
$zongji.Range("D3").Value = 0.01

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 3
$zongji.Range("D2").Value = 0.02
